$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "NO"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "expected outcome"
$ws.Range("E1").Value = "Models"

# --- Row 2 ---
$ws.Range("C2").Value = "Check if all necessary input files exists"
$ws.Range("D2").Value = "Breaks with logfile message like this:`n25-Jul-2017 13:36:56 Start input checks`n25-Jul-2017 13:36:56 ERROR: ""Children_OralSingle_IV_Multi_withTypo.xml"" does not exist`n25-Jul-2017 13:36:56 ERROR: ""Children_OralSingle_IV_Multi_withTypo.csv "" does not exist"
$ws.Range("E2").Value = "TestExample"

# --- Row 3 ---
$ws.Range("C3").Value = "Find inconsistencies in output definition, unkonw output path, unknown units, or wrong unit for dimensions"
$ws.Range("D3").Value = "Breaks with logfile message`nERROR: Outputpath ""Organism|PeripheralVenousBlood|Hydroxy_Itraconazole|Plasma (Peripheral Venous Blood) with Typo"" could not be found in model`nERROR: For unit ""µmol/l"", there is no common dimension with display unit ""cm""`nERROR: unit ""typo"" for seems to be no default OSPSuite unit`nERROR: unit ""typo2"" for seems to be no default OSPSuite unit`nERROR: For unit ""µmol/l"", there is no common dimension with display unit ""h"""
$ws.Range("E3").Value = "TestExample"

# --- Row 4 ---
$ws.Range("C4").Value = "Population simulation and PK Parameter calculation of a population with single application"
$ws.Range("D4").Value = "Two csv files are generated: `nsimulations/SingleIvBolus-Results.csv`nsimulations/SingleIvBolus-PK-Analyses.csv`nautomatic test compares the result to PK-SIM export of this model. It works for simulations, fails for PK Parameter with extrapolation ???"
$ws.Range("C4").Font.Color = 255

# --- Row 5 ---
$ws.Range("C5").Value = "Population simulation and PK Parameter calculation of a multi application"
$ws.Range("D5").Value = "Two csv files are generated: `nSimulations/OralSingle_IV_Multi-PK-Analyses.csv`nSimulations/OralSingle_IV_Multi-Results.csv`nautomatic test compares the result to PK-SIM export of this model. It works for simulations, fails for PK Parameter with extrapolation ???"
$ws.Range("E5").Value = "TestExample"

# --- Row 6 ---
$ws.Range("C6").Value = "Populations simulation with studydesign.csv which sets the dose according to the dose per bodyweigth"
$ws.Range("D6").Value = "Two csv files are generated: `nsimulations/SingleIvBolus-Results.csv`nsimulations/SingleIvBolus-PK-Analyses.csv`nAutomatic test checks if the simulation results is a factor 10 highe the the PK-Sim export"
$ws.Range("C6").Font.Color = 255

# --- Row 7 ---
$ws.Range("C7").Value = "Populations simulation with studydesign.csv dose per surface area"
$ws.Range("D7").Value = "Two csv files are generated: `nsimulations/SingleIvBolus-Results.csv`nsimulations/SingleIvBolus-PK-Analyses.csv`nAutomatic test checks if the simulation results is a factor 10 highe the the PK-Sim export"
$ws.Range("E7").Value = "7.2_BSA_Example"
$ws.Range("C7").Font.Color = 255

# --- Row 8 ---
$ws.Range("C8").Value = "Read nonmen file with individual timeprofiles simulate mean model and plot data vs prediction"
$ws.Range("D8").Value = "data is converted to a mat file: tmp/PO320mg/dataTP.mat  check for random individuals, if time and dv is correctly transferred`nfigures are created: time profile, same figures as in corresponding output, `n"

# --- Row 9 ---
$ws.Range("C9").Value = "Population VPC with data for a single population"
$ws.Range("D9").Value = "figures are created: pyhsiological, time profile and pkParameter, same figures as in corresponding output, `n"
$ws.Range("C9").Font.Color = 255

# --- Row 10 ---
$ws.Range("C10").Value = "Population simulation for workflowtype parallelComparison"
$ws.Range("D10").Value = "figures are created: pyhsiological, time profile and pkParameter, same figures as in corresponding output, `n"

# --- Row 11 ---
$ws.Range("C11").Value = "Population simulation for workflowtype pediatric"
$ws.Range("D11").Value = "figures are created: pyhsiological, time profile and pkParameter, same figures as in corresponding output, `n"

# --- Row 12 (clear) ---
$ws.Range("C12").Value = ""

# --- Window / selection ---
$excel.ActiveWindow.Left = 192
$ws.Range("D10").Select()
